$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Fill in stat values for row 39 (Flynn / Puppeteer)
$ws.Range("G39").Value = 45
$ws.Range("H39").Value = 60
$ws.Range("I39").Value = 75
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 45
$ws.Range("L39").Value = 60
$ws.Range("M39").Value = 50

# Fill in stat values for row 41 (Lis / Hunter)
$ws.Range("G41").Value = 70
$ws.Range("H41").Value = 50
$ws.Range("I41").Value = 60
$ws.Range("J41").Value = 45
$ws.Range("K41").Value = 55
$ws.Range("L41").Value = 40
$ws.Range("M41").Value = 70

# Check the "has art" checkboxes for rows 39 and 41
$ws.CheckBoxes("Check Box 75").Value = 1
$ws.CheckBoxes("Check Box 77").Value = 1

# Update the sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("L35").Select()
